$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Smart FORTWO']"
$arr[0,1] = "[0.2]"
$arr[0,2] = "[0.8000000000000002]"
$arr[0,3] = "[10.560000000000002]"
$arr[0,4] = [double]10.56
$ws.Range("B7:F7").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Tesla MODEL 3', 'Fiat 500 E']"
$arr[0,1] = "[0.15, 0.2]"
$arr[0,2] = "[0.9500000000000003, 0.9000000000000002]"
$arr[0,3] = "[40.000000000000014, 16.800000000000004]"
$arr[0,4] = [double]56.80000000000002
$ws.Range("B8:F8").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Fiat 500 E', 'Hyundai KONA 39 kWh', 'Fiat 500 E']"
$arr[0,1] = "[0.3, 0.1, 0.25]"
$arr[0,2] = "[1.0, 0.9000000000000002, 0.9500000000000003]"
$arr[0,3] = "[16.799999999999997, 31.20000000000001, 16.800000000000008]"
$arr[0,4] = [double]64.80000000000001
$ws.Range("B9:F9").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.15]"
$arr[0,2] = "[0.8000000000000002]"
$arr[0,3] = "[34.20986111111112]"
$arr[0,4] = [double]34.20986111111112
$ws.Range("B10:F10").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'VW E-UP', 'VW ID.3']"
$arr[0,1] = "[0.35, 0.2, 0.25]"
$arr[0,2] = "[0.7500000000000001, 1.0, 1.0]"
$arr[0,3] = "[21.052222222222227, 29.439999999999998, 43.5]"
$arr[0,4] = [double]93.99222222222222
$ws.Range("B11:F11").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Hyundai KONA 39 kWh', 'KIA EV6']"
$arr[0,1] = "[0.25, 0.15, 0.15]"
$arr[0,2] = "[0.65, 0.8500000000000002, 0.8000000000000002]"
$arr[0,3] = "[21.052222222222223, 27.300000000000008, 48.10000000000001]"
$arr[0,4] = [double]96.45222222222225
$ws.Range("B12:F12").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.2]"
$arr[0,2] = "[0.8500000000000002]"
$arr[0,3] = "[34.20986111111112]"
$arr[0,4] = [double]34.20986111111112
$ws.Range("B13:F13").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B17:F17").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'SKODA ENYAQ 77kWh', 'Polestar 2']"
$arr[0,1] = "[0.25, 0.25, 0.35]"
$arr[0,2] = "[1.0, 1.0, 0.8500000000000002]"
$arr[0,3] = "[39.47291666666666, 57.75, 37.500000000000014]"
$arr[0,4] = [double]134.7229166666667
$ws.Range("B18:F18").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['VW ID.3', 'Hyundai KONA 64 kWh', 'Hyundai IONIQ5 77kWh']"
$arr[0,1] = "[0.3, 0.1, 0.25]"
$arr[0,2] = "[1.0, 0.9500000000000003, 0.9500000000000003]"
$arr[0,3] = "[40.599999999999994, 54.40000000000002, 53.90000000000002]"
$arr[0,4] = [double]148.9
$ws.Range("B19:F19").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['VW ID.4']"
$arr[0,1] = "[0.25]"
$arr[0,2] = "[0.8000000000000002]"
$arr[0,3] = "[42.35000000000001]"
$arr[0,4] = [double]42.35000000000001
$ws.Range("B20:F20").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B21:F21").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Tesla MODEL 3', 'Others']"
$arr[0,1] = "[0.2, 0.3]"
$arr[0,2] = "[0.9000000000000002, 0.9500000000000003]"
$arr[0,3] = "[35.00000000000001, 34.20986111111113]"
$arr[0,4] = [double]69.20986111111114
$ws.Range("B31:F31").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Renault ZOE', 'Fiat 500 E', 'Audi Q4']"
$arr[0,1] = "[0.2, 0.05, 0.05]"
$arr[0,2] = "[1.0, 0.8000000000000002, 0.9000000000000002]"
$arr[0,3] = "[41.6, 18.000000000000004, 65.11000000000001]"
$arr[0,4] = [double]124.71
$ws.Range("B32:F32").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Audi Q4', 'Others', 'VW ID.5']"
$arr[0,1] = "[0.05, 0.1, 0.1]"
$arr[0,2] = "[0.9000000000000002, 0.8500000000000002, 0.7000000000000001]"
$arr[0,3] = "[65.11000000000001, 39.47291666666668, 46.20000000000001]"
$arr[0,4] = [double]150.7829166666667
$ws.Range("B33:F33").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B34:F34").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B35:F35").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.25]"
$arr[0,2] = "[0.7500000000000001]"
$arr[0,3] = "[26.315277777777784]"
$arr[0,4] = [double]26.31527777777778
$ws.Range("B36:F36").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Opel MOKKA', 'Tesla MODEL 3', 'Renault ZOE', 'Others', 'Audi E-TRON']"
$arr[0,1] = "[0.25, 0.2, 0.25, 0.2, 0.35]"
$arr[0,2] = "[0.7000000000000001, 0.9000000000000002, 0.9500000000000003, 0.8500000000000002, 0.7500000000000001]"
$arr[0,3] = "[20.250000000000004, 35.00000000000001, 36.40000000000001, 34.20986111111112, 34.000000000000014]"
$arr[0,4] = [double]159.8598611111112
$ws.Range("B37:F37").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.4]"
$arr[0,2] = "[0.9000000000000002]"
$arr[0,3] = "[26.315277777777787]"
$arr[0,4] = [double]26.31527777777779
$ws.Range("B39:F39").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B41:F41").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Others']"
$arr[0,1] = "[0.35, 0.3]"
$arr[0,2] = "[0.8500000000000002, 1.0]"
$arr[0,3] = "[26.315277777777787, 36.841388888888886]"
$arr[0,4] = [double]63.15666666666667
$ws.Range("B42:F42").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Hyundai KONA 64 kWh', 'TESLA MODEL Y', 'Audi E-TRON', 'Fiat 500 E']"
$arr[0,1] = "[0.15, 0.05, 0.25, 0.25]"
$arr[0,2] = "[0.8500000000000002, 1.0, 0.9000000000000002, 0.7000000000000001]"
$arr[0,3] = "[44.80000000000001, 71.25, 55.25000000000002, 10.8]"
$arr[0,4] = [double]182.1000000000001
$ws.Range("B43:F43").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Renault ZOE']"
$arr[0,1] = "[0.4]"
$arr[0,2] = "[0.9000000000000002]"
$arr[0,3] = "[26.00000000000001]"
$arr[0,4] = [double]26.00000000000001
$ws.Range("B44:F44").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['MINI Cooper SE', 'Hyundai KONA 39 kWh']"
$arr[0,1] = "[0.1, 0.15]"
$arr[0,2] = "[0.7000000000000001, 0.8500000000000002]"
$arr[0,3] = "[17.340000000000003, 27.300000000000008]"
$arr[0,4] = [double]44.64000000000001
$ws.Range("B55:F55").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.15]"
$arr[0,2] = "[0.7500000000000001]"
$arr[0,3] = "[31.578333333333337]"
$arr[0,4] = [double]31.57833333333334
$ws.Range("B56:F56").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Hyundai KONA 64 kWh', 'Others']"
$arr[0,1] = "[0.15, 0.35, 0.05]"
$arr[0,2] = "[0.8500000000000002, 0.8500000000000002, 0.7500000000000001]"
$arr[0,3] = "[36.84138888888889, 32.000000000000014, 36.84138888888889]"
$arr[0,4] = [double]105.6827777777778
$ws.Range("B57:F57").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Renault ZOE']"
$arr[0,1] = "[0.35]"
$arr[0,2] = "[0.8500000000000002]"
$arr[0,3] = "[26.00000000000001]"
$arr[0,4] = [double]26.00000000000001
$ws.Range("B58:F58").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'MINI Cooper SE']"
$arr[0,1] = "[0.15, 0.25]"
$arr[0,2] = "[0.9000000000000002, 0.8000000000000002]"
$arr[0,3] = "[39.47291666666668, 15.895000000000003]"
$arr[0,4] = [double]55.36791666666668
$ws.Range("B60:F60").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Polestar 2', 'Audi Q4', 'Others', 'VW ID.4', 'KIA EV6']"
$arr[0,1] = "[0.15, 0.4, 0.2, 0.15, 0.3]"
$arr[0,2] = "[1.0, 0.7000000000000001, 0.8500000000000002, 0.7000000000000001, 0.8500000000000002]"
$arr[0,3] = "[63.75, 22.98, 34.20986111111112, 42.35, 40.70000000000002]"
$arr[0,4] = [double]203.9898611111111
$ws.Range("B61:F61").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B62:F62").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Smart FORTWO']"
$arr[0,1] = "[0.45]"
$arr[0,2] = "[0.8500000000000002]"
$arr[0,3] = "[7.040000000000004]"
$arr[0,4] = [double]7.040000000000004
$ws.Range("B64:F64").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "[]"
$arr[0,1] = "[]"
$arr[0,2] = "[]"
$arr[0,3] = "[]"
$arr[0,4] = [double]0
$ws.Range("B65:F65").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Others']"
$arr[0,1] = "[0.1, 0.2]"
$arr[0,2] = "[0.8000000000000002, 0.8000000000000002]"
$arr[0,3] = "[36.84138888888889, 31.578333333333337]"
$arr[0,4] = [double]68.41972222222223
$ws.Range("B66:F66").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Others', 'Fiat 500 E']"
$arr[0,1] = "[0.35, 0.4999999999999999, 0.25]"
$arr[0,2] = "[0.9500000000000003, 1.0, 0.8500000000000002]"
$arr[0,3] = "[31.578333333333347, 26.315277777777784, 14.400000000000006]"
$arr[0,4] = [double]72.29361111111113
$ws.Range("B67:F67").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.15]"
$arr[0,2] = "[0.8500000000000002]"
$arr[0,3] = "[36.84138888888889]"
$arr[0,4] = [double]36.84138888888889
$ws.Range("B68:F68").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['TESLA MODEL Y', 'Hyundai KONA 39 kWh']"
$arr[0,1] = "[0.2, 0.15]"
$arr[0,2] = "[0.9000000000000002, 0.9000000000000002]"
$arr[0,3] = "[52.500000000000014, 29.250000000000007]"
$arr[0,4] = [double]81.75000000000003
$ws.Range("B78:F78").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['VW ID.3']"
$arr[0,1] = "[0.25]"
$arr[0,2] = "[0.9000000000000002]"
$arr[0,3] = "[37.70000000000002]"
$arr[0,4] = [double]37.70000000000002
$ws.Range("B79:F79").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Others']"
$arr[0,1] = "[0.2, 0.1]"
$arr[0,2] = "[0.8000000000000002, 1.0]"
$arr[0,3] = "[31.578333333333337, 47.3675]"
$arr[0,4] = [double]78.94583333333334
$ws.Range("B80:F80").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Dacia SPRING', 'Others']"
$arr[0,1] = "[0.15, 0.15]"
$arr[0,2] = "[0.9000000000000002, 0.9000000000000002]"
$arr[0,3] = "[20.100000000000005, 39.47291666666668]"
$arr[0,4] = [double]59.57291666666669
$ws.Range("B81:F81").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Audi E-TRON']"
$arr[0,1] = "[0.2]"
$arr[0,2] = "[0.6]"
$arr[0,3] = "[34.0]"
$arr[0,4] = [double]34
$ws.Range("B83:F83").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Hyundai KONA 64 kWh', 'VW E-UP', 'Others', 'TESLA MODEL Y']"
$arr[0,1] = "[0.25, 0.2, 0.35, 0.25]"
$arr[0,2] = "[0.9000000000000002, 0.8500000000000002, 0.8500000000000002, 0.9500000000000003]"
$arr[0,3] = "[41.600000000000016, 23.92, 26.315277777777787, 52.50000000000002]"
$arr[0,4] = [double]144.3352777777778
$ws.Range("B84:F84").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.1]"
$arr[0,2] = "[0.9500000000000003]"
$arr[0,3] = "[44.73597222222224]"
$arr[0,4] = [double]44.73597222222224
$ws.Range("B85:F85").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['VW ID.3']"
$arr[0,1] = "[0.25]"
$arr[0,2] = "[0.9500000000000003]"
$arr[0,3] = "[40.600000000000016]"
$arr[0,4] = [double]40.60000000000002
$ws.Range("B86:F86").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.45]"
$arr[0,2] = "[0.9500000000000003]"
$arr[0,3] = "[26.315277777777787]"
$arr[0,4] = [double]26.31527777777779
$ws.Range("B89:F89").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others', 'Audi Q4', 'Others']"
$arr[0,1] = "[0.35, 0.25, 0.3]"
$arr[0,2] = "[0.65, 0.9500000000000003, 0.7500000000000001]"
$arr[0,3] = "[15.789166666666668, 53.62000000000002, 23.683750000000007]"
$arr[0,4] = [double]93.0929166666667
$ws.Range("B90:F90").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['VW ID.5', 'Others']"
$arr[0,1] = "[0.2, 0.25]"
$arr[0,2] = "[0.65, 1.0]"
$arr[0,3] = "[34.65, 39.47291666666666]"
$arr[0,4] = [double]74.12291666666667
$ws.Range("B91:F91").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "['Others']"
$arr[0,1] = "[0.4]"
$arr[0,2] = "[0.8500000000000002]"
$arr[0,3] = "[23.683750000000007]"
$arr[0,4] = [double]23.68375000000001
$ws.Range("B92:F92").Value = $arr
